$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 337, shifting existing rows 337:358 down to 338:359
$ws.Rows.Item(337).Insert()

# Populate the newly inserted row 337 with the new data record
$ws.Cells.Item(337, 1).Value = 3
$ws.Cells.Item(337, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(337, 3).Value = "Coquimbo"
$ws.Cells.Item(337, 4).Value = 44826
$ws.Cells.Item(337, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(337, 5).Value = 5
$ws.Cells.Item(337, 6).Value = 100112001
$ws.Cells.Item(337, 7).Value = "Berenjena"
$ws.Cells.Item(337, 8).Value = "Sin especificar"
$ws.Cells.Item(337, 9).Value = "Primera"
$ws.Cells.Item(337, 10).Value = 45
$ws.Cells.Item(337, 11).Value = 13000
$ws.Cells.Item(337, 12).Value = 13000
$ws.Cells.Item(337, 13).Value = 13000
$ws.Cells.Item(337, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(337, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(337, 16).Value = 217
$ws.Cells.Item(337, 17).Value = 60
$ws.Cells.Item(337, 18).Value = "Hortaliza"
